$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 1252258.9
$ws.Range("L46").Value = 4292430.300000001
$ws.Range("N46").Value = -4292668.300000001
$ws.Range("J46").Value = 1430810.1

$ws.Range("J60").Value = 1430810.1
$ws.Range("H60").Value = 1252258.9
$ws.Range("L60").Value = 4292430.300000001
$ws.Range("N60").Value = -4293398.300000001

$ws.Range("H87").Value = 33354
$ws.Range("L87").Value = 33354
$ws.Range("N87").Value = -35850
$ws.Range("J87").Value = 33354

$ws.Range("H90").Value = 33354
$ws.Range("N90").Value = -112542
$ws.Range("L90").Value = 100062
$ws.Range("J90").Value = 33354

$ws.Range("H92").Value = 995.38464
$ws.Range("I92").Value = 994.5454999999999
$ws.Range("K92").Value = 994.5454999999999
$ws.Range("M92").Value = 253.4545000000001

$ws.Range("I137").Value = 2037.2307
$ws.Range("H137").Value = 60704.883
$ws.Range("K137").Value = 6111.6921
$ws.Range("M137").Value = -3561.6921

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J32").Value = 181482.17
$ws.Range("H32").Value = 38557.12
$ws.Range("M32").Value = -21455.412
$ws.Range("K32").Value = 21742.412
$ws.Range("L32").Value = 181482.17
$ws.Range("N32").Value = -182056.17
$ws.Range("I32").Value = 21742.412

$ws.Range("N45").Value = -3133.6667
$ws.Range("K45").Value = 1687028.1
$ws.Range("H45").Value = 676239.0600000001
$ws.Range("L45").Value = 2379.6667
$ws.Range("J45").Value = 2379.6667
$ws.Range("M45").Value = -1686651.1
$ws.Range("I45").Value = 1687028.1

$ws.Range("H61").Value = 653.4194
$ws.Range("K61").Value = 614.37933
$ws.Range("M61").Value = -402.37933
$ws.Range("I61").Value = 614.37933

$ws.Range("I74").Value = 868.04346
$ws.Range("K74").Value = 868.04346
$ws.Range("H74").Value = 1904.9688
$ws.Range("M74").Value = 5.956540000000018

$ws.Range("H77").Value = 1904.9688
$ws.Range("I77").Value = 868.04346
$ws.Range("K77").Value = 4340.2173
$ws.Range("M77").Value = 27.78269999999975

$ws.Range("H92").Value = 98000
$ws.Range("L92").Value = 98000
$ws.Range("J92").Value = 98000
$ws.Range("N92").Value = -102992

$ws.Range("L102").Value = 1854
$ws.Range("N102").Value = -5098
$ws.Range("J102").Value = 1854
$ws.Range("H102").Value = 2071.25

$ws.Range("H132").Value = 10092.174
$ws.Range("K132").Value = 34803.315
$ws.Range("I132").Value = 11601.105
$ws.Range("M132").Value = -32273.315
$ws.Range("N132").Value = -13834.25
$ws.Range("J132").Value = 2924.75
$ws.Range("L132").Value = 8774.25

$ws.Range("H136").Value = 653.4194
$ws.Range("M136").Value = 706.8620100000001
$ws.Range("I136").Value = 614.37933
$ws.Range("K136").Value = 1843.13799

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N107").Value = -5078.4667
$ws.Range("L107").Value = 1238.4667
$ws.Range("J107").Value = 1238.4667
$ws.Range("H107").Value = 19158.125

$ws.Range("M134").Value = -2110.200000000001
$ws.Range("L134").Value = 2394
$ws.Range("K134").Value = 4645.200000000001
$ws.Range("N134").Value = -7464
$ws.Range("J134").Value = 798
$ws.Range("H134").Value = 1423.3334
$ws.Range("I134").Value = 1548.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1201.9565
$ws.Range("I31").Value = 1131.3636
$ws.Range("K31").Value = 1131.3636
$ws.Range("M31").Value = -836.3635999999999

$ws.Range("I34").Value = 1131.3636
$ws.Range("M34").Value = -929.3635999999999
$ws.Range("K34").Value = 1131.3636
$ws.Range("H34").Value = 1201.9565

$ws.Range("J64").Value = 25271
$ws.Range("N64").Value = -25767
$ws.Range("H64").Value = 25271
$ws.Range("L64").Value = 25271

$ws.Range("N67").Value = -26987
$ws.Range("J67").Value = 25271
$ws.Range("H67").Value = 25271
$ws.Range("L67").Value = 25271

$ws.Range("L86").Value = 8499.75
$ws.Range("I86").Value = 6292.1665
$ws.Range("H86").Value = 7175.2
$ws.Range("J86").Value = 8499.75
$ws.Range("N86").Value = -10745.75
$ws.Range("M86").Value = -5169.1665

$ws.Range("N89").Value = -53730.75
$ws.Range("J89").Value = 8499.75
$ws.Range("L89").Value = 42498.75
$ws.Range("M89").Value = -25844.8325
$ws.Range("H89").Value = 7175.2
$ws.Range("I89").Value = 6292.1665
$ws.Range("K89").Value = 31460.8325

$ws.Range("H132").Value = 3951.1904
$ws.Range("K132").Value = 12587.5269
$ws.Range("I132").Value = 4195.8423
$ws.Range("M132").Value = -10057.5269
$ws.Range("N132").Value = -9941
$ws.Range("J132").Value = 1627
$ws.Range("L132").Value = 4881

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I3").Value = 1359.4
$ws.Range("K3").Value = 4078.2
$ws.Range("H3").Value = 1359.4
$ws.Range("M3").Value = -3966.2

$ws.Range("L103").Value = 333364350
$ws.Range("H103").Value = 47626204
$ws.Range("N103").Value = -333366108
$ws.Range("J103").Value = 111121450
$ws.Range("K103").Value = 14306.25
$ws.Range("M103").Value = -13427.25
$ws.Range("I103").Value = 4768.75

$ws.Range("M117").Value = -1356.9998
$ws.Range("K117").Value = 4798.9998
$ws.Range("N117").Value = -25598.693
$ws.Range("I117").Value = 1599.6666
$ws.Range("L117").Value = 18714.693
$ws.Range("H117").Value = 5368.5
$ws.Range("J117").Value = 6238.231

$ws.Range("H136").Value = 3963.9375
$ws.Range("M136").Value = -943.2497999999996
$ws.Range("I136").Value = 2014.4166
$ws.Range("K136").Value = 6043.2498

$ws.Range("N137").Value = -24792.7998
$ws.Range("I137").Value = 9092166
$ws.Range("H137").Value = 3849491.8
$ws.Range("J137").Value = 4864.2666
$ws.Range("K137").Value = 27276498
$ws.Range("L137").Value = 14592.7998
$ws.Range("M137").Value = -27271398

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M43").Value = -1515.3334
$ws.Range("I43").Value = 1666.3334
$ws.Range("K43").Value = 1666.3334
$ws.Range("H43").Value = 20999.4

$ws.Range("H132").Value = 2228.3872
$ws.Range("K132").Value = 6732.5172
$ws.Range("I132").Value = 2244.1724
$ws.Range("M132").Value = -4202.5172
$ws.Range("N132").Value = -11058.5
$ws.Range("J132").Value = 1999.5
$ws.Range("L132").Value = 5998.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K16").Value = 472.2857
$ws.Range("M16").Value = -302.2857
$ws.Range("I16").Value = 472.2857
$ws.Range("H16").Value = 508.55554

$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H82").Value = 3483.3333

$ws.Range("H85").Value = 3483.3333

$ws.Range("H132").Value = 3890.2727
$ws.Range("K132").Value = 7891.5
$ws.Range("I132").Value = 2630.5
$ws.Range("M132").Value = -5361.5
$ws.Range("N132").Value = -26808.9995
$ws.Range("J132").Value = 7249.6665
$ws.Range("L132").Value = 21748.9995

$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("H138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N45").Value = -27580.6
$ws.Range("H45").Value = 29665.5
$ws.Range("L45").Value = 26598.6
$ws.Range("J45").Value = 26598.6

$ws.Range("H54").Value = 43332.668
$ws.Range("J54").Value = 49999
$ws.Range("L54").Value = 49999
$ws.Range("N54").Value = -51039

$ws.Range("J131").Value = 87452.5
$ws.Range("N131").Value = -97532.5
$ws.Range("H131").Value = 87452.5
$ws.Range("L131").Value = 87452.5

$ws.Range("H132").Value = 10620.082
$ws.Range("K132").Value = 36766.287
$ws.Range("I132").Value = 12255.429
$ws.Range("M132").Value = -34236.287
